$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item("Semaine_1")

# --- Extend the "Semaine_1" table by 4 rows (84-87) for the new deliveries ---
$tbl.ListRows.Add() | Out-Null
$tbl.ListRows.Add() | Out-Null
$tbl.ListRows.Add() | Out-Null
$tbl.ListRows.Add() | Out-Null

# Row 84 : Mame Mareme NDIAYE - DKR Plateau - Babacar Mbaye Kébé
$ws.Cells.Item(84, 1).Value = 45954
$ws.Cells.Item(84, 2).Value = "Mame Mareme NDIAYE"
$ws.Cells.Item(84, 3).Value = "DKR PLATEAU"
$ws.Cells.Item(84, 4).Value = "DKR Plateau"
$ws.Cells.Item(84, 5).Value = "Babacar Mbaye Kébé"
$ws.Cells.Item(84, 6).Value = 776169696
$ws.Cells.Item(84, 7).Value = "Grossiste"
$ws.Cells.Item(84, 8).Value = "Client Partenaire"
$ws.Cells.Item(84, 9).Value = "Livraison"
$ws.Cells.Item(84, 10).Value = "Livraison"
$ws.Cells.Item(84, 11).Value = "Lait Janus 18gx100"
$ws.Cells.Item(84, 12).Value = 25
$ws.Cells.Item(84, 13).Value = 6500
$ws.Cells.Item(84, 14).Value = 162500
$ws.Cells.Item(84, 15).Formula = '="S"&_xlfn.ISOWEEKNUM([@Date])'
$ws.Cells.Item(84, 16).Formula = '=TEXT([@Date],"MMMM")'

# Row 85 : Fatoumata TRAORE - Malika - Abdou sow 1
$ws.Cells.Item(85, 1).Value = 45954
$ws.Cells.Item(85, 2).Value = "Fatoumata TRAORE"
$ws.Cells.Item(85, 3).Value = "KEUR MASSAR"
$ws.Cells.Item(85, 4).Value = "Malika"
$ws.Cells.Item(85, 5).Value = "Abdou sow 1"
$ws.Cells.Item(85, 6).Value = 779646150
$ws.Cells.Item(85, 7).Value = "Grossiste"
$ws.Cells.Item(85, 8).Value = "Client Partenaire"
$ws.Cells.Item(85, 9).Value = "Livraison"
$ws.Cells.Item(85, 10).Value = "RAS "
$ws.Cells.Item(85, 11).Value = "Café stick Refraish 1,5gx09boites"
$ws.Cells.Item(85, 12).Value = 25
$ws.Cells.Item(85, 13).Value = 26000
$ws.Cells.Item(85, 14).Value = 650000
$ws.Cells.Item(85, 15).Formula = '="S"&_xlfn.ISOWEEKNUM([@Date])'
$ws.Cells.Item(85, 16).Formula = '=TEXT([@Date],"MMMM")'

# Row 86 : Fatoumata TRAORE - Malika - DJIBRIL laye
$ws.Cells.Item(86, 1).Value = 45954
$ws.Cells.Item(86, 2).Value = "Fatoumata TRAORE"
$ws.Cells.Item(86, 3).Value = "KEUR MASSAR"
$ws.Cells.Item(86, 4).Value = "Malika"
$ws.Cells.Item(86, 5).Value = "DJIBRIL laye"
$ws.Cells.Item(86, 6).Value = 778657940
$ws.Cells.Item(86, 7).Value = "Grossiste"
$ws.Cells.Item(86, 8).Value = "Client Partenaire"
$ws.Cells.Item(86, 9).Value = "Livraison"
$ws.Cells.Item(86, 10).Value = "Ras"
$ws.Cells.Item(86, 11).Value = "Café stick Refraish 1,5gx09boites"
$ws.Cells.Item(86, 12).Value = 25
$ws.Cells.Item(86, 13).Value = 26000
$ws.Cells.Item(86, 14).Value = 650000
$ws.Cells.Item(86, 15).Formula = '="S"&_xlfn.ISOWEEKNUM([@Date])'
$ws.Cells.Item(86, 16).Formula = '=TEXT([@Date],"MMMM")'

# Row 87 : Ngouye NDIAYE - Zac Mbao - DJIBRIL Traoré
$ws.Cells.Item(87, 1).Value = 45954
$ws.Cells.Item(87, 2).Value = "Ngouye NDIAYE"
$ws.Cells.Item(87, 3).Value = "PIKINE"
$ws.Cells.Item(87, 4).Value = "Zac Mbao"
$ws.Cells.Item(87, 5).Value = "DJIBRIL Traoré"
$ws.Cells.Item(87, 6).Value = 775262371
$ws.Cells.Item(87, 7).Value = "Grossiste"
$ws.Cells.Item(87, 8).Value = "Client Partenaire"
$ws.Cells.Item(87, 9).Value = "Livraison"
$ws.Cells.Item(87, 10).Value = "Ok"
$ws.Cells.Item(87, 11).Value = "Café stick Refraish 1,5gx09boites"
$ws.Cells.Item(87, 12).Value = 10
$ws.Cells.Item(87, 13).Value = 26000
$ws.Cells.Item(87, 14).Value = 260000
$ws.Cells.Item(87, 15).Formula = '="S"&_xlfn.ISOWEEKNUM([@Date])'
$ws.Cells.Item(87, 16).Formula = '=TEXT([@Date],"MMMM")'

# --- Stray note typed below the table ---
$ws.Cells.Item(95, 10).Value = " "

# --- Update the view: scroll down to the new rows and select the note cell ---
$ws.Range("J96").Select() | Out-Null
